$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Style-Cell($rng, [int]$colorIndex) {
    $rng.WrapText = $true
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.Font.ColorIndex = $colorIndex
}

# Row 9 - TC_03
$rng = $ws.Range("A9")
$rng.Value = "TC_03"
Style-Cell $rng 17
$rng = $ws.Range("B9")
$rng.Value = "Verify the error message when passing the multiple(2 attributes) invalid attribute names"
Style-Cell $rng 17
$rng = $ws.Range("C9")
$rng.Value = "GraphQL"
Style-Cell $rng 17
$rng = $ws.Range("D9")
Style-Cell $rng 17
$rng = $ws.Range("E9")
Style-Cell $rng 17
$rng = $ws.Range("F9")
Style-Cell $rng 17
$rng = $ws.Range("G9")
Style-Cell $rng 17
$rng = $ws.Range("H9")
$rng.Value = "'200"
Style-Cell $rng 17
$rng = $ws.Range("I9")
$rng.Value = "`n{`n`t`"meta`":`n`t{`n`t`t`"version`":`"1.0.0`",`n`t`t`"errors`":`n`t`t[`n`t`t`t{`n`t`t`t`t`"timestamp`":`"2020-02-04T13:39:45.663Z`n`t`t`t`t[`n`t`t`t`t`tGMT`n`t`t`t`t]`",`n`t`t`t`t`"error`":`"ValidationError`",`n`t`t`t`t`"message`":`"Validation error of type FieldUndefined: Field 'geopoliticalRelationshipTypeCode1' in type 'GeoPoliticalRelationshipType' is undefined @ 'relationshipTypes/geopoliticalRelationshipTypeCode1'`",`n`t`t`t`t`"path`":null`n`t`t`t},`n`t`t`t{`n`t`t`t`t`"timestamp`":`"2020-02-04T13:39:45.663Z`n`t`t`t`t[`n`t`t`t`t`tGMT`n`t`t`t`t]`",`n`t`t`t`t`"error`":`"ValidationError`",`n`t`t`t`t`"message`":`"Validation error of type FieldUndefined: Field 'areaRelationshipTypeDescription1' in type 'GeoPoliticalRelationshipType' is undefined @ 'relationshipTypes/areaRelationshipTypeDescription1'`",`n`t`t`t`t`"path`":null`n`t`t`t}`n`t`t]`n`t},`n`t`"data`":null`n}"
Style-Cell $rng 17
$rng = $ws.Range("J9")
$rng.Value = "Fail"
Style-Cell $rng 17
$rng = $ws.Range("K9")
Style-Cell $rng 17

# Row 10 - TC_04
$rng = $ws.Range("A10")
$rng.Value = "TC_04"
Style-Cell $rng 10
$rng = $ws.Range("B10")
$rng.Value = "Verify no results fetched when passing the invalid geopoliticalRelationshipTypeCode  parameter"
Style-Cell $rng 10
$rng = $ws.Range("C10")
$rng.Value = "GraphQL"
Style-Cell $rng 10
$rng = $ws.Range("D10")
$rng.Value = "`n{`n`t`"query`":`"`n`t{`n`t`t  relationshipTypes (relationshipTypeCode :\`"ABC\`") `n`t`t{`n`t`t`t    geopoliticalRelationshipTypeCode    areaRelationshipTypeDescription  `n`t`t}`n`t}`"`n}"
Style-Cell $rng 10
$rng = $ws.Range("E10")
$rng.Value = "Input_GeoRsTypeCode: ABC`nInput_GeoRsTypeDesc: null`n"
Style-Cell $rng 10
$rng = $ws.Range("F10")
$rng.Value = "NA"
Style-Cell $rng 10
$rng = $ws.Range("G10")
Style-Cell $rng 10
$rng = $ws.Range("H10")
$rng.Value = "'200"
Style-Cell $rng 10
$rng = $ws.Range("I10")
$rng.Value = "`n{`n`t`"meta`":`n`t{`n`t`t`"version`":`"1.0.0`",`n`t`t`"errors`":`n`t`t[`n`t`t]`n`t},`n`t`"data`":`n`t{`n`t`t`"relationshipTypes`":`n`t`t[`n`t`t]`n`t}`n}"
Style-Cell $rng 10
$rng = $ws.Range("J10")
$rng.Value = "Pass"
Style-Cell $rng 10
$rng = $ws.Range("K10")
Style-Cell $rng 10

Write-Host "done"